$d = $word.ActiveDocument

# The new "Help Center Resources" block to append after the customer-support
# hotline bullet and before the trailing blank paragraphs.
$texts = @(
    "Help Center Resources:",
    "Getting Started Guide",
    "A comprehensive guide on navigating our website, creating an account, and placing orders.",
    "Payment and Checkout Guide",
    "Detailed information on accepted payment methods, secure checkout procedures, and payment-related FAQs.",
    "Shipping and Delivery Information",
    "An overview of our shipping process, delivery times, and international shipping options.",
    "Returns and Refunds Policies",
    "Detailed policies on returns, exchanges, and refunds, including step-by-step instructions for initiating a return.",
    "Account Management",
    "Resources on managing your account settings, updating personal information, and changing passwords.",
    "Security and Privacy",
    "Information on how we secure your personal and payment data to ensure a safe and private shopping experience.",
    "Product FAQs",
    "Specific FAQs related to product categories, sizing guides, and other relevant details."
)

# Locate the anchor paragraph: "Calling our customer support hotline at [phone number]."
# (the last bulleted item under "How do I contact customer support...").
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "Calling our customer support hotline*") {
        $anchorIndex = $i
    }
}

if ($anchorIndex -eq -1) {
    throw "Could not find anchor paragraph 'Calling our customer support hotline...'"
}

$cur = $d.Paragraphs.Item($anchorIndex)

for ($i = 0; $i -lt $texts.Length; $i++) {
    $cur.Range.InsertParagraphAfter()
    $cur = $d.Paragraphs.Item($anchorIndex + 1 + $i)
    # Reset to the default (Normal) paragraph style so the new paragraph does
    # not inherit the "ListParagraph" style/bullet of the anchor paragraph.
    $cur.Style = "Normal"
    $cur.Range.Text = $texts[$i]
}
